$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so numeric-looking strings (e.g. "210.50", "0.650")
# are preserved exactly instead of being normalized into numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '81.718.98'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +5.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.198.68'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.00%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.50'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '637.19'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.293'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +28.81%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.191.62'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.595'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +11.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000266'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +19.50%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.780.46'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.28'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +5.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '81.466.69'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +5.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.192.84'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.25'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +14.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.45'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.31'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '443.81'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.26'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +10.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.11'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.11'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +9.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.32'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.361.28'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '77.31'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000129'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +12.98%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.26'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.90%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '576.45'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +11.11%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.62%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.05'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.88%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.153'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +12.90%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Cronos'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.141'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +31.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.37'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.46%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.418'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +6.69%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.14'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +25.47%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.08'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +18.99%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.98'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +11.50%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '160.18'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.99%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '190.01'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '45.43'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +6.39%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +6.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.786'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.650'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.65%  '
